# Update the "Metadata" sheet of the CodeSystem-destino-traslado workbook
# to the new publication snapshot:
#   - Status:         draft -> active
#   - Date:            2024-12-13T10:10:51-03:00 -> 2024-12-16T14:50:05-03:00
#   - Case Sensitive:  false -> true
# (new version "iq and urgency")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B6").Value = "active"
$ws.Range("B8").Value = "2024-12-16T14:50:05-03:00"
# Leading apostrophe forces this to be stored as literal text "true"
# (otherwise Excel auto-types the bare word true/false as a Boolean).
$ws.Range("B17").Value = "'true"
